$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Cells.Item(32, 2).Value = 6149867
$ws.Cells.Item(32, 5).Value = 'Gangwon FC'
$ws.Cells.Item(32, 6).Value = 'Gwangju FC'
$ws.Cells.Item(32, 8).Value = 1
$ws.Cells.Item(32, 9).Value = 'D'
$ws.Cells.Item(32, 10).Value = 2.75
$ws.Cells.Item(32, 11).Value = 3.1
$ws.Cells.Item(32, 12).Value = 2.5
$ws.Cells.Item(32, 13).Value = 3.3
$ws.Cells.Item(32, 14).Value = 3.1
$ws.Cells.Item(32, 15).Value = 2.2
$ws.Cells.Item(32, 16).Value = 0.25
$ws.Cells.Item(32, 17).Value = 1.95
$ws.Cells.Item(32, 18).Value = 1.9
$ws.Cells.Item(32, 19).Value = 2.25
$ws.Cells.Item(32, 20).Value = 2
$ws.Cells.Item(32, 21).Value = 1.85
$ws.Cells.Item(32, 23).Value = 2.1
$ws.Cells.Item(32, 24).Value = -1
$ws.Cells.Item(32, 25).Value = 0.475
$ws.Cells.Item(32, 26).Value = -0.5
$ws.Cells.Item(32, 27).Value = -0.5
$ws.Cells.Item(32, 28).Value = 0.425

# Row 33
$ws.Cells.Item(33, 2).Value = 6149414
$ws.Cells.Item(33, 5).Value = 'Jeju United'
$ws.Cells.Item(33, 6).Value = 'Daegu FC'
$ws.Cells.Item(33, 8).Value = 2
$ws.Cells.Item(33, 9).Value = 'A'
$ws.Cells.Item(33, 10).Value = 2.1
$ws.Cells.Item(33, 11).Value = 3.4
$ws.Cells.Item(33, 12).Value = 3.2
$ws.Cells.Item(33, 13).Value = 2.3
$ws.Cells.Item(33, 14).Value = 3.3
$ws.Cells.Item(33, 15).Value = 3
$ws.Cells.Item(33, 16).Value = -0.25
$ws.Cells.Item(33, 17).Value = 2.05
$ws.Cells.Item(33, 18).Value = 1.8
$ws.Cells.Item(33, 19).Value = 2.5
$ws.Cells.Item(33, 20).Value = 1.9
$ws.Cells.Item(33, 21).Value = 1.95
$ws.Cells.Item(33, 23).Value = -1
$ws.Cells.Item(33, 24).Value = 2
$ws.Cells.Item(33, 25).Value = -1
$ws.Cells.Item(33, 26).Value = 0.8
$ws.Cells.Item(33, 27).Value = 0.8999999999999999
$ws.Cells.Item(33, 28).Value = -1

# Row 98
$ws.Cells.Item(98, 2).Value = 6353260
$ws.Cells.Item(98, 5).Value = 'Pohang Steelers'
$ws.Cells.Item(98, 6).Value = 'Ulsan Hyundai'
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 2.375
$ws.Cells.Item(98, 11).Value = 3.5
$ws.Cells.Item(98, 12).Value = 2.8
$ws.Cells.Item(98, 13).Value = 2.55
$ws.Cells.Item(98, 14).Value = 3.4
$ws.Cells.Item(98, 15).Value = 2.625
$ws.Cells.Item(98, 16).Value = 0
$ws.Cells.Item(98, 17).Value = 1.825
$ws.Cells.Item(98, 18).Value = 2.025
$ws.Cells.Item(98, 19).Value = 2.5
$ws.Cells.Item(98, 20).Value = 2.025
$ws.Cells.Item(98, 21).Value = 1.825
$ws.Cells.Item(98, 23).Value = 2.4
$ws.Cells.Item(98, 25).Value = 0
$ws.Cells.Item(98, 26).Value = 0
$ws.Cells.Item(98, 28).Value = 0.825

# Row 99
$ws.Cells.Item(99, 2).Value = 6353261
$ws.Cells.Item(99, 5).Value = 'Suwon FC'
$ws.Cells.Item(99, 6).Value = 'FC Seoul'
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 1
$ws.Cells.Item(99, 10).Value = 3.4
$ws.Cells.Item(99, 11).Value = 3.6
$ws.Cells.Item(99, 12).Value = 2
$ws.Cells.Item(99, 13).Value = 3.75
$ws.Cells.Item(99, 14).Value = 3.8
$ws.Cells.Item(99, 15).Value = 1.85
$ws.Cells.Item(99, 16).Value = 0.5
$ws.Cells.Item(99, 17).Value = 2.025
$ws.Cells.Item(99, 18).Value = 1.825
$ws.Cells.Item(99, 19).Value = 2.75
$ws.Cells.Item(99, 20).Value = 1.825
$ws.Cells.Item(99, 21).Value = 2.025
$ws.Cells.Item(99, 23).Value = 2.8
$ws.Cells.Item(99, 25).Value = 1.025
$ws.Cells.Item(99, 26).Value = -1
$ws.Cells.Item(99, 28).Value = 1.025

# Row 105
$ws.Cells.Item(105, 2).Value = 6384418
$ws.Cells.Item(105, 5).Value = 'Gwangju FC'
$ws.Cells.Item(105, 6).Value = 'Gangwon FC'
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 'H'
$ws.Cells.Item(105, 10).Value = 1.909
$ws.Cells.Item(105, 11).Value = 3.4
$ws.Cells.Item(105, 12).Value = 4
$ws.Cells.Item(105, 13).Value = 1.95
$ws.Cells.Item(105, 14).Value = 3.4
$ws.Cells.Item(105, 15).Value = 4
$ws.Cells.Item(105, 16).Value = -0.5
$ws.Cells.Item(105, 17).Value = 2
$ws.Cells.Item(105, 18).Value = 1.85
$ws.Cells.Item(105, 19).Value = 2.25
$ws.Cells.Item(105, 20).Value = 1.925
$ws.Cells.Item(105, 21).Value = 1.925
$ws.Cells.Item(105, 22).Value = 0.95
$ws.Cells.Item(105, 24).Value = -1
$ws.Cells.Item(105, 25).Value = 1
$ws.Cells.Item(105, 26).Value = -1
$ws.Cells.Item(105, 28).Value = 0.925

# Row 106
$ws.Cells.Item(106, 2).Value = 6387791
$ws.Cells.Item(106, 5).Value = 'FC Seoul'
$ws.Cells.Item(106, 6).Value = 'Jeonbuk Motors'
$ws.Cells.Item(106, 7).Value = 0
$ws.Cells.Item(106, 8).Value = 2
$ws.Cells.Item(106, 9).Value = 'A'
$ws.Cells.Item(106, 10).Value = 2.375
$ws.Cells.Item(106, 11).Value = 3.3
$ws.Cells.Item(106, 12).Value = 2.875
$ws.Cells.Item(106, 13).Value = 2.4
$ws.Cells.Item(106, 14).Value = 3.3
$ws.Cells.Item(106, 15).Value = 2.8
$ws.Cells.Item(106, 16).Value = 0
$ws.Cells.Item(106, 17).Value = 1.775
$ws.Cells.Item(106, 18).Value = 2.1
$ws.Cells.Item(106, 19).Value = 2.5
$ws.Cells.Item(106, 20).Value = 1.85
$ws.Cells.Item(106, 21).Value = 2
$ws.Cells.Item(106, 22).Value = -1
$ws.Cells.Item(106, 24).Value = 1.8
$ws.Cells.Item(106, 25).Value = -1
$ws.Cells.Item(106, 26).Value = 1.1
$ws.Cells.Item(106, 28).Value = 1

# Row 107
$ws.Cells.Item(107, 2).Value = 6384127
$ws.Cells.Item(107, 5).Value = 'Daejeon Hana Citizen'
$ws.Cells.Item(107, 6).Value = 'Jeju United'
$ws.Cells.Item(107, 10).Value = 2.375
$ws.Cells.Item(107, 12).Value = 2.875
$ws.Cells.Item(107, 13).Value = 2.6
$ws.Cells.Item(107, 15).Value = 2.625
$ws.Cells.Item(107, 16).Value = 0
$ws.Cells.Item(107, 17).Value = 1.9
$ws.Cells.Item(107, 18).Value = 1.95
$ws.Cells.Item(107, 19).Value = 2.75
$ws.Cells.Item(107, 20).Value = 1.975
$ws.Cells.Item(107, 21).Value = 1.875
$ws.Cells.Item(107, 22).Value = 1.6
$ws.Cells.Item(107, 25).Value = 0.8999999999999999
$ws.Cells.Item(107, 28).Value = 0.875

# Row 108
$ws.Cells.Item(108, 2).Value = 6384126
$ws.Cells.Item(108, 5).Value = 'Suwon Bluewings'
$ws.Cells.Item(108, 6).Value = 'Pohang Steelers'
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = 'H'
$ws.Cells.Item(108, 10).Value = 4
$ws.Cells.Item(108, 11).Value = 3.4
$ws.Cells.Item(108, 12).Value = 1.909
$ws.Cells.Item(108, 13).Value = 4
$ws.Cells.Item(108, 14).Value = 3.3
$ws.Cells.Item(108, 15).Value = 1.95
$ws.Cells.Item(108, 16).Value = 0.5
$ws.Cells.Item(108, 17).Value = 1.875
$ws.Cells.Item(108, 18).Value = 1.975
$ws.Cells.Item(108, 19).Value = 2.5
$ws.Cells.Item(108, 20).Value = 2.1
$ws.Cells.Item(108, 21).Value = 1.775
$ws.Cells.Item(108, 22).Value = 3
$ws.Cells.Item(108, 23).Value = -1
$ws.Cells.Item(108, 25).Value = 0.875
$ws.Cells.Item(108, 26).Value = -1
$ws.Cells.Item(108, 27).Value = -1
$ws.Cells.Item(108, 28).Value = 0.7749999999999999

# Row 109
$ws.Cells.Item(109, 2).Value = 6384125
$ws.Cells.Item(109, 5).Value = 'Daegu FC'
$ws.Cells.Item(109, 6).Value = 'Suwon FC'
$ws.Cells.Item(109, 7).Value = 2
$ws.Cells.Item(109, 8).Value = 2
$ws.Cells.Item(109, 9).Value = 'D'
$ws.Cells.Item(109, 10).Value = 1.666
$ws.Cells.Item(109, 11).Value = 3.75
$ws.Cells.Item(109, 12).Value = 4.75
$ws.Cells.Item(109, 13).Value = 1.6
$ws.Cells.Item(109, 14).Value = 3.8
$ws.Cells.Item(109, 15).Value = 5.25
$ws.Cells.Item(109, 16).Value = -0.75
$ws.Cells.Item(109, 17).Value = 1.85
$ws.Cells.Item(109, 18).Value = 2
$ws.Cells.Item(109, 20).Value = 2.025
$ws.Cells.Item(109, 21).Value = 1.825
$ws.Cells.Item(109, 22).Value = -1
$ws.Cells.Item(109, 23).Value = 2.8
$ws.Cells.Item(109, 25).Value = -1
$ws.Cells.Item(109, 26).Value = 1
$ws.Cells.Item(109, 27).Value = 1.025
$ws.Cells.Item(109, 28).Value = -1

# Row 125
$ws.Cells.Item(125, 2).Value = 7333492
$ws.Cells.Item(125, 5).Value = 'Incheon Utd'
$ws.Cells.Item(125, 6).Value = 'Jeonbuk Motors'
$ws.Cells.Item(125, 7).Value = 1
$ws.Cells.Item(125, 8).Value = 1
$ws.Cells.Item(125, 9).Value = 'D'
$ws.Cells.Item(125, 10).Value = 2.75
$ws.Cells.Item(125, 11).Value = 3.2
$ws.Cells.Item(125, 12).Value = 2.6
$ws.Cells.Item(125, 13).Value = 3.2
$ws.Cells.Item(125, 14).Value = 3.25
$ws.Cells.Item(125, 15).Value = 2.25
$ws.Cells.Item(125, 16).Value = 0.25
$ws.Cells.Item(125, 17).Value = 1.85
$ws.Cells.Item(125, 18).Value = 2
$ws.Cells.Item(125, 19).Value = 2.25
$ws.Cells.Item(125, 20).Value = 1.825
$ws.Cells.Item(125, 21).Value = 2.025
$ws.Cells.Item(125, 23).Value = 2.25
$ws.Cells.Item(125, 24).Value = -1
$ws.Cells.Item(125, 25).Value = 0.425
$ws.Cells.Item(125, 26).Value = -0.5
$ws.Cells.Item(125, 27).Value = -0.5
$ws.Cells.Item(125, 28).Value = 0.5125

# Row 126
$ws.Cells.Item(126, 2).Value = 7334084
$ws.Cells.Item(126, 5).Value = 'Suwon FC'
$ws.Cells.Item(126, 6).Value = 'Suwon Bluewings'
$ws.Cells.Item(126, 7).Value = 2
$ws.Cells.Item(126, 8).Value = 3
$ws.Cells.Item(126, 9).Value = 'A'
$ws.Cells.Item(126, 10).Value = 2.15
$ws.Cells.Item(126, 11).Value = 3.5
$ws.Cells.Item(126, 12).Value = 3.2
$ws.Cells.Item(126, 13).Value = 2.55
$ws.Cells.Item(126, 14).Value = 3.4
$ws.Cells.Item(126, 15).Value = 2.625
$ws.Cells.Item(126, 16).Value = 0
$ws.Cells.Item(126, 17).Value = 1.875
$ws.Cells.Item(126, 18).Value = 1.975
$ws.Cells.Item(126, 19).Value = 2.75
$ws.Cells.Item(126, 20).Value = 1.9
$ws.Cells.Item(126, 21).Value = 1.95
$ws.Cells.Item(126, 23).Value = -1
$ws.Cells.Item(126, 24).Value = 1.625
$ws.Cells.Item(126, 25).Value = -1
$ws.Cells.Item(126, 26).Value = 0.9750000000000001
$ws.Cells.Item(126, 27).Value = 0.8999999999999999
$ws.Cells.Item(126, 28).Value = -1

# Row 131
$ws.Cells.Item(131, 2).Value = 7333496
$ws.Cells.Item(131, 5).Value = 'Pohang Steelers'
$ws.Cells.Item(131, 6).Value = 'Daegu FC'
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 10).Value = 1.85
$ws.Cells.Item(131, 11).Value = 3.4
$ws.Cells.Item(131, 12).Value = 3.6
$ws.Cells.Item(131, 13).Value = 2.1
$ws.Cells.Item(131, 16).Value = -0.25
$ws.Cells.Item(131, 17).Value = 1.8
$ws.Cells.Item(131, 18).Value = 2.05
$ws.Cells.Item(131, 20).Value = 1.975
$ws.Cells.Item(131, 21).Value = 1.875
$ws.Cells.Item(131, 22).Value = 1.1
$ws.Cells.Item(131, 25).Value = 0.8
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = 0.875

# Row 133
$ws.Cells.Item(133, 2).Value = 7334086
$ws.Cells.Item(133, 5).Value = 'Gangwon FC'
$ws.Cells.Item(133, 6).Value = 'Suwon FC'
$ws.Cells.Item(133, 7).Value = 2
$ws.Cells.Item(133, 10).Value = 2.05
$ws.Cells.Item(133, 11).Value = 3.6
$ws.Cells.Item(133, 12).Value = 3.3
$ws.Cells.Item(133, 13).Value = 2.05
$ws.Cells.Item(133, 16).Value = -0.5
$ws.Cells.Item(133, 17).Value = 2.05
$ws.Cells.Item(133, 18).Value = 1.8
$ws.Cells.Item(133, 20).Value = 1.825
$ws.Cells.Item(133, 21).Value = 2.025
$ws.Cells.Item(133, 22).Value = 1.05
$ws.Cells.Item(133, 25).Value = 1.05
$ws.Cells.Item(133, 27).Value = -0.5
$ws.Cells.Item(133, 28).Value = 0.5125

# Row 153
$ws.Cells.Item(153, 2).Value = 7716466
$ws.Cells.Item(153, 5).Value = 'FC Seoul'
$ws.Cells.Item(153, 6).Value = 'Jeju United'
$ws.Cells.Item(153, 7).Value = 2
$ws.Cells.Item(153, 8).Value = 0
$ws.Cells.Item(153, 9).Value = 'H'
$ws.Cells.Item(153, 10).Value = 2.1
$ws.Cells.Item(153, 11).Value = 3.3
$ws.Cells.Item(153, 12).Value = 3.2
$ws.Cells.Item(153, 13).Value = 2.2
$ws.Cells.Item(153, 15).Value = 3.1
$ws.Cells.Item(153, 16).Value = -0.25
$ws.Cells.Item(153, 17).Value = 1.975
$ws.Cells.Item(153, 18).Value = 1.875
$ws.Cells.Item(153, 19).Value = 2.25
$ws.Cells.Item(153, 20).Value = 1.85
$ws.Cells.Item(153, 21).Value = 2
$ws.Cells.Item(153, 22).Value = 1.2
$ws.Cells.Item(153, 23).Value = -1
$ws.Cells.Item(153, 25).Value = 0.9750000000000001
$ws.Cells.Item(153, 26).Value = -1
$ws.Cells.Item(153, 27).Value = -0.5
$ws.Cells.Item(153, 28).Value = 0.5

# Row 154
$ws.Cells.Item(154, 2).Value = 7716465
$ws.Cells.Item(154, 5).Value = 'Daegu FC'
$ws.Cells.Item(154, 6).Value = 'Suwon FC'
$ws.Cells.Item(154, 7).Value = 1
$ws.Cells.Item(154, 8).Value = 1
$ws.Cells.Item(154, 9).Value = 'D'
$ws.Cells.Item(154, 10).Value = 2
$ws.Cells.Item(154, 11).Value = 3.25
$ws.Cells.Item(154, 12).Value = 3.5
$ws.Cells.Item(154, 13).Value = 2
$ws.Cells.Item(154, 15).Value = 3.5
$ws.Cells.Item(154, 16).Value = -0.5
$ws.Cells.Item(154, 17).Value = 2.05
$ws.Cells.Item(154, 18).Value = 1.8
$ws.Cells.Item(154, 19).Value = 2.5
$ws.Cells.Item(154, 20).Value = 1.975
$ws.Cells.Item(154, 21).Value = 1.875
$ws.Cells.Item(154, 22).Value = -1
$ws.Cells.Item(154, 23).Value = 2.3
$ws.Cells.Item(154, 25).Value = -1
$ws.Cells.Item(154, 26).Value = 0.8
$ws.Cells.Item(154, 27).Value = -1
$ws.Cells.Item(154, 28).Value = 0.875

# Row 155
$ws.Cells.Item(155, 2).Value = 7715267
$ws.Cells.Item(155, 5).Value = 'Gimcheon Sangmu FC'
$ws.Cells.Item(155, 6).Value = 'Jeonbuk Motors'
$ws.Cells.Item(155, 10).Value = 3.5
$ws.Cells.Item(155, 11).Value = 3.25
$ws.Cells.Item(155, 12).Value = 2
$ws.Cells.Item(155, 13).Value = 3.25
$ws.Cells.Item(155, 14).Value = 3.25
$ws.Cells.Item(155, 15).Value = 2.1
$ws.Cells.Item(155, 16).Value = 0.25
$ws.Cells.Item(155, 17).Value = 2.05
$ws.Cells.Item(155, 18).Value = 1.8
$ws.Cells.Item(155, 19).Value = 2.5
$ws.Cells.Item(155, 20).Value = 2.05
$ws.Cells.Item(155, 21).Value = 1.8
$ws.Cells.Item(155, 22).Value = 2.25
$ws.Cells.Item(155, 25).Value = 1.05
$ws.Cells.Item(155, 28).Value = 0.8

# Row 156
$ws.Cells.Item(156, 2).Value = 7715264
$ws.Cells.Item(156, 5).Value = 'Pohang Steelers'
$ws.Cells.Item(156, 6).Value = 'Gwangju FC'
$ws.Cells.Item(156, 10).Value = 2.3
$ws.Cells.Item(156, 11).Value = 3.1
$ws.Cells.Item(156, 12).Value = 3
$ws.Cells.Item(156, 13).Value = 2.625
$ws.Cells.Item(156, 14).Value = 3.1
$ws.Cells.Item(156, 15).Value = 2.6
$ws.Cells.Item(156, 16).Value = 0
$ws.Cells.Item(156, 17).Value = 1.975
$ws.Cells.Item(156, 18).Value = 1.875
$ws.Cells.Item(156, 19).Value = 2.25
$ws.Cells.Item(156, 20).Value = 1.975
$ws.Cells.Item(156, 21).Value = 1.875
$ws.Cells.Item(156, 22).Value = 1.625
$ws.Cells.Item(156, 25).Value = 0.9750000000000001
$ws.Cells.Item(156, 28).Value = 0.875

# Row 176
$ws.Cells.Item(176, 2).Value = 7715280
$ws.Cells.Item(176, 5).Value = 'Jeonbuk Motors'
$ws.Cells.Item(176, 6).Value = 'Gwangju FC'
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 'H'
$ws.Cells.Item(176, 10).Value = 2.05
$ws.Cells.Item(176, 12).Value = 3.6
$ws.Cells.Item(176, 13).Value = 1.85
$ws.Cells.Item(176, 14).Value = 3.3
$ws.Cells.Item(176, 15).Value = 4.2
$ws.Cells.Item(176, 16).Value = -0.5
$ws.Cells.Item(176, 17).Value = 1.9
$ws.Cells.Item(176, 18).Value = 1.95
$ws.Cells.Item(176, 20).Value = 1.95
$ws.Cells.Item(176, 21).Value = 1.9
$ws.Cells.Item(176, 22).Value = 0.8500000000000001
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 0.8999999999999999
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.95

# Row 177
$ws.Cells.Item(177, 2).Value = 7715281
$ws.Cells.Item(177, 5).Value = 'FC Seoul'
$ws.Cells.Item(177, 6).Value = 'Pohang Steelers'
$ws.Cells.Item(177, 8).Value = 4
$ws.Cells.Item(177, 9).Value = 'A'
$ws.Cells.Item(177, 10).Value = 2.375
$ws.Cells.Item(177, 12).Value = 2.9
$ws.Cells.Item(177, 13).Value = 2.75
$ws.Cells.Item(177, 14).Value = 3
$ws.Cells.Item(177, 15).Value = 2.75
$ws.Cells.Item(177, 16).Value = 0
$ws.Cells.Item(177, 17).Value = 1.925
$ws.Cells.Item(177, 18).Value = 1.925
$ws.Cells.Item(177, 20).Value = 1.85
$ws.Cells.Item(177, 21).Value = 2
$ws.Cells.Item(177, 22).Value = -1
$ws.Cells.Item(177, 24).Value = 1.75
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = 0.925
$ws.Cells.Item(177, 27).Value = 0.8500000000000001

# Row 214
$ws.Cells.Item(214, 2).Value = 7716534
$ws.Cells.Item(214, 5).Value = 'FC Seoul'
$ws.Cells.Item(214, 6).Value = 'Daegu FC'
$ws.Cells.Item(214, 7).Value = 1
$ws.Cells.Item(214, 8).Value = 2
$ws.Cells.Item(214, 10).Value = 2.1
$ws.Cells.Item(214, 12).Value = 3.6
$ws.Cells.Item(214, 13).Value = 2.05
$ws.Cells.Item(214, 14).Value = 3.2
$ws.Cells.Item(214, 15).Value = 3.9
$ws.Cells.Item(214, 19).Value = 2.25
$ws.Cells.Item(214, 20).Value = 1.95
$ws.Cells.Item(214, 21).Value = 1.9
$ws.Cells.Item(214, 24).Value = 2.9
$ws.Cells.Item(214, 27).Value = 0.95
$ws.Cells.Item(214, 28).Value = -1

# Row 215
$ws.Cells.Item(215, 2).Value = 7715307
$ws.Cells.Item(215, 5).Value = 'Gangwon FC'
$ws.Cells.Item(215, 6).Value = 'Ulsan Hyundai'
$ws.Cells.Item(215, 8).Value = 0
$ws.Cells.Item(215, 9).Value = 'H'
$ws.Cells.Item(215, 10).Value = 3.3
$ws.Cells.Item(215, 11).Value = 3.4
$ws.Cells.Item(215, 12).Value = 2.15
$ws.Cells.Item(215, 13).Value = 3.6
$ws.Cells.Item(215, 14).Value = 3.8
$ws.Cells.Item(215, 15).Value = 1.909
$ws.Cells.Item(215, 16).Value = 0.5
$ws.Cells.Item(215, 17).Value = 1.9
$ws.Cells.Item(215, 18).Value = 1.95
$ws.Cells.Item(215, 19).Value = 2.75
$ws.Cells.Item(215, 20).Value = 1.85
$ws.Cells.Item(215, 21).Value = 2
$ws.Cells.Item(215, 22).Value = 2.6
$ws.Cells.Item(215, 24).Value = -1
$ws.Cells.Item(215, 25).Value = 0.8999999999999999
$ws.Cells.Item(215, 26).Value = -1
$ws.Cells.Item(215, 27).Value = -1
$ws.Cells.Item(215, 28).Value = 1

# Row 216
$ws.Cells.Item(216, 2).Value = 7715306
$ws.Cells.Item(216, 5).Value = 'Gwangju FC'
$ws.Cells.Item(216, 6).Value = 'Jeonbuk Motors'
$ws.Cells.Item(216, 7).Value = 0
$ws.Cells.Item(216, 8).Value = 3
$ws.Cells.Item(216, 9).Value = 'A'
$ws.Cells.Item(216, 10).Value = 2.45
$ws.Cells.Item(216, 11).Value = 3.25
$ws.Cells.Item(216, 12).Value = 2.9
$ws.Cells.Item(216, 13).Value = 1.95
$ws.Cells.Item(216, 14).Value = 3.6
$ws.Cells.Item(216, 15).Value = 3.6
$ws.Cells.Item(216, 16).Value = -0.5
$ws.Cells.Item(216, 17).Value = 2
$ws.Cells.Item(216, 18).Value = 1.85
$ws.Cells.Item(216, 20).Value = 1.825
$ws.Cells.Item(216, 21).Value = 2.025
$ws.Cells.Item(216, 22).Value = -1
$ws.Cells.Item(216, 24).Value = 2.6
$ws.Cells.Item(216, 25).Value = -1
$ws.Cells.Item(216, 26).Value = 0.8500000000000001
$ws.Cells.Item(216, 27).Value = 0.4125
$ws.Cells.Item(216, 28).Value = -0.5

# Row 223 standalone odds update
$ws.Cells.Item(223, 17).Value = 1.875
$ws.Cells.Item(223, 18).Value = 1.975